$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 13010
$ws1.Range("F5").Value = 84
$ws1.Range("F6").Value = 94
$ws1.Range("F10").Value = 12984
$ws1.Range("F13").Value = 8721
$ws1.Range("F14").Value = 7737
$ws1.Range("F15").Value = 206
$ws1.Range("F16").Value = 117
$ws1.Range("F19").Value = 990
$ws1.Range("F24").Value = 332

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 13010
$ws4.Range("F6").Value = 84
$ws4.Range("F7").Value = 94
$ws4.Range("F11").Value = 12984
$ws4.Range("F14").Value = 8721
$ws4.Range("F15").Value = 7737
$ws4.Range("F16").Value = 206
$ws4.Range("F17").Value = 117
$ws4.Range("F20").Value = 990
$ws4.Range("F27").Value = 332
